$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newE = @{
    2  = 2.1
    3  = 2.5
    4  = 2.2
    5  = 2.3
    6  = 2
    7  = 2.4
    8  = 2.2
    9  = 2.6
    10 = 2.1
    11 = 2.7
    12 = 2.91
    13 = 2.23
    14 = 3.03
    15 = 2.88
    16 = 3.08
    17 = 2.88
    18 = 2.73
    19 = 2.93
    20 = 3.14
    21 = 2.88
    22 = 3.37
    23 = 2.87
    24 = 2.86
    25 = 2.99
    26 = 3.32
    27 = 2.77
    28 = 2.69
    29 = 3.31
    30 = 3.12
    31 = 2.67
}

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 4).Value = "粉土"
    $ws.Cells.Item($row, 5).Value = $newE[$row]
}
